# doc/开发进度.xlsx -- progress tracker update
# "add properties for ERect in QStudioSCADA and QSCADARunTime."
#
# Row 8 = 矩形 (ERect/"Rectangle"): status moves 进行中 (in progress) -> 已完成 (done)
# Row 9 = 图片 (Picture): status moves 未开始 (not started) -> 进行中 (in progress)
# Selection moves from B7 to B8 to reflect the row just finished.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("画面编辑器")

# 矩形 (ERect) is now complete.
$ws.Range("B8").Value = "已完成"
$ws.Range("B8").Interior.Color = 5287936   # RGB(0,176,80) green - matches "已完成" fill

# 图片 now moves into progress.
$ws.Range("B9").Value = "进行中"
$ws.Range("B9").Interior.Color = 65535     # RGB(255,255,0) yellow - matches "进行中" fill

# Reflect the new active cell/selection.
$ws.Range("B8").Select()
